$d = $word.ActiveDocument

# --- Change 1: paragraph "1024px X 500px PN" + "G" -> merge trailing two runs into "px PNG" ---
# (keeps the leading "1024px X 500" run untouched, matches target run layout)
$p5 = $d.Paragraphs(5)
$p5Start = $p5.Range.Start
$p5End = $p5.Range.End
$p5Range = $d.Range($p5Start, $p5End - 1)
$xml5 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1024px X 500</w:t></w:r><w:r><w:t>px PNG</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p5Range.InsertXML($xml5)

# --- Change 2: "Welcome Image" paragraph gains " Play Store" run and the _GoBack bookmark moves here ---
$p6 = $d.Paragraphs(6)
$p6Start = $p6.Range.Start
$p6End = $p6.Range.End
$p6Range = $d.Range($p6Start, $p6End - 1)
$xml6 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Welcome Image</w:t></w:r><w:r><w:t xml:space="preserve"> Play Store</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p6Range.InsertXML($xml6)

# --- Change 3: remove the old _GoBack bookmark that used to sit at the end of the last paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
